# The sheet's lone header row had a stray leading blank column: real data
# lived in B:F while A only held a duplicate of F (GENE) with header-style
# bleed-through on rows 2-3. Fix the "MODEL_CONDITION" header typo and drop
# the spurious column A so the data lines up in A:E.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the mis-typed header text (still sitting in column E at this point).
$ws.Cells.Item(1, 5).Value = "MODELCONDITION"

# Remove column A entirely - this shifts B:F left into A:E, carrying the
# header styling (bold/border) from B1:F1 into A1:E1, and the old column F
# (which duplicated column A's numbers) lands in the new column E.
$ws.Columns(1).Delete()
